$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67

# Row 6
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("X6").Value = 5.5
$ws.Range("AD6").Value = 10
$ws.Range("AG6").Value = 26
$ws.Range("AI6").Value = 34
$ws.Range("AQ6").Value = 15
$ws.Range("AW6").Value = 11
$ws.Range("AZ6").Value = 301

# Row 10
$ws.Range("I10").Value = 3.1
$ws.Range("M10").Value = 1.11
$ws.Range("N10").Value = 6.5
$ws.Range("U10").Value = 2.1
$ws.Range("V10").Value = 1.67
$ws.Range("AC10").Value = 6.5
$ws.Range("AI10").Value = 12
$ws.Range("AK10").Value = 29

# Row 15
$ws.Range("O15").Value = 1.5
$ws.Range("P15").Value = 2.5
$ws.Range("BD15").Value = 126

# Row 17
$ws.Range("O17").Value = 1.33
$ws.Range("P17").Value = 3.25
